$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column Q, row 4: new year header (2020) ---
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2020

# --- Column Q, row 5: Education value ---
$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 53.2

# --- Column Q, row 6: Health value ---
$ws.Range("P6").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = 23.2

# --- Column Q, row 7: Social protection value (needs "0.0" number format,
#     matching the newly introduced cellXfs entry) ---
$ws.Range("P7").Copy($ws.Range("Q7"))
$ws.Range("Q7").Value = 10
$ws.Range("Q7").NumberFormat = "0.0"

# --- Column Q, row 8: bottom row, thick-bottom border variant (also needs
#     the "0.0" number format) ---
$ws.Range("P8").Copy($ws.Range("Q8"))
$ws.Range("Q8").Value = 20
$ws.Range("Q8").NumberFormat = "0.0"

# --- Selection moves from P5 to P9 ---
$ws.Range("P9").Select() | Out-Null
